$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = 0.1293707326694715
$ws.Range("E2").Value = 9.078375525823422
$ws.Range("F2").Value = 24.01755214409346
